$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the header row (row 1) ---------------------------------------
# Columns A:J were "<Name>_old" -> "<Name>_FV2410"
# Column K ("diff") is unchanged
# Columns L:U were "<Name>_new" -> "<Name>_FV2504"
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")

foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    $txt = [string]$cell.Value2
    $cell.Value = $txt.Replace("_old", "_FV2410")
}

foreach ($col in $newCols) {
    $cell = $ws.Range($col + "1")
    $txt = [string]$cell.Value2
    $cell.Value = $txt.Replace("_new", "_FV2504")
}

# --- Turn the used range into an Excel Table ------------------------------
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U65"), $null, 1)
$lo.Name = "Table1"

# --- Freeze the header row -------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
